# pharma_pos.product_size demo data update
# - "Bottle" is renamed to "Tablet" for the "500 ml" size row
# - A new "20 cc" / "Injection" row is appended, dated 2021-01-31 (44227)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 ("500 ml") Type column: Bottle -> Tablet
$ws.Range("B2").Value = "Tablet"

# Row 3 ("10 mg") Type column stays "Tablet" (kept explicit for clarity)
$ws.Range("B3").Value = "Tablet"

# New row 4: copy formatting (style + number format) from row 3, then set values
$ws.Range("A3:C3").Copy()
$ws.Range("A4:C4").PasteSpecial(-4122)

$ws.Range("A4").Value = "20 cc"
$ws.Range("B4").Value = "Injection"
$ws.Range("C4").Value = 44227
